$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") - entire row shifts everything up by one.
$ws.Range("A26").EntireRow.Delete()

# After the above delete, the row that was "SC 92" (originally row 28) is now row 27.
# Delete it too, shifting everything below up by one more.
$ws.Range("A27").EntireRow.Delete()

# --- Column F value changes (independent of the row shifts above) ---
$ws.Range("F2").Value = 18.03
$ws.Range("F6").Value = $null
$ws.Range("F12").Value = 17.45
$ws.Range("F14").Value = $null
$ws.Range("F20").Value = 17.73
$ws.Range("F21").Value = 16.58
$ws.Range("F23").Value = $null
$ws.Range("F24").Value = $null
$ws.Range("F31").Value = 17.18
$ws.Range("F33").Value = 17.53

# --- Column D value changes (independent of the row shifts above) ---
$ws.Range("D26").Value = -13.8
$ws.Range("D27").Value = $null
$ws.Range("D30").Value = -13.6
$ws.Range("D32").Value = $null
